# Duplicate row 2's populated cells into row 3 (preserving per-cell
# value types), then overwrite the handful of cells whose content
# actually differs for the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","P","Q","R","S","T","U","V","W", `
          "Y","Z","AA","AB","AC","AD","AE","AG","AT","AW","AX","AY")
foreach ($col in $cols) {
  $ws.Range($col + "2").Copy($ws.Range($col + "3"))
}

$ws.Range("A3").Value = 112181620
$ws.Range("B3").Value = 78578
$ws.Range("Q3").Value = 818894.527582898
$ws.Range("R3").Value = 7382401.525517201
$ws.Range("Y3").Value = "'2023-08-08"
$ws.Range("AA3").Value = "'2023-08-08"
$ws.Range("AC3").Value = "Påträffad under Sveaskogs naturvärdesinventering"
